# Update the "War Room" Dataframe sheet:
#  - Turn the header dates (previously text like "5 de junio") into real
#    date values formatted as d-mmm, and extend the header with two more
#    label columns ("9-jun" / "12-jun").
#  - Refresh the daily figures in columns E:F and add two brand-new days
#    of data in columns G:H for every product row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
# A1/B1 keep their text ("Línea" / "Descripción") but pick up an explicit
# text number format.
$ws.Range("A1:B1").NumberFormat = "@"

# C1:F1 become real dates (5-8 June 2023) instead of shared-string labels.
$ws.Range("C1").Value = 45082
$ws.Range("D1").Value = 45083
$ws.Range("E1").Value = 45084
$ws.Range("F1").Value = 45085
$ws.Range("C1:F1").NumberFormat = "d-mmm"

# G1:H1 are brand-new header cells for the two additional days.
$ws.Range("G1").Value = "9-jun"
$ws.Range("H1").Value = "12-jun"
$ws.Range("G1:H1").NumberFormat = "@"
$ws.Range("G1:H1").HorizontalAlignment = -4152

# --- New data columns G:H (match the existing C:F number style) ------
$ws.Range("G2:H18").NumberFormat = "0"
$ws.Range("G2:H18").HorizontalAlignment = -4108

# --- Row 2 ------------------------------------------------------------
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 0

# --- Row 3 ------------------------------------------------------------
$ws.Range("E3").Value = 12.660506184546698
$ws.Range("F3").Value = 12.304801759214161
$ws.Range("G3").Value = 13.448627194409992
$ws.Range("H3").Value = 15.261967779093684

# --- Row 4 ------------------------------------------------------------
$ws.Range("E4").Value = 17.268694184539548
$ws.Range("F4").Value = 17.990243152180277
$ws.Range("G4").Value = 18.505905110729802
$ws.Range("H4").Value = 18.152062957099719

# --- Row 5 ------------------------------------------------------------
$ws.Range("E5").Value = 20.284107897115465
$ws.Range("F5").Value = 19.282784015060884
$ws.Range("G5").Value = 18.877445359546506
$ws.Range("H5").Value = 21.289386798917953

# --- Row 6 ------------------------------------------------------------
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0

# --- Row 7 ------------------------------------------------------------
$ws.Range("E7").Value = 9.5953065940943123
$ws.Range("F7").Value = 10.707621182446195
$ws.Range("G7").Value = 9.1969667110451123
$ws.Range("H7").Value = 11.353553897619783

# --- Row 8 ------------------------------------------------------------
$ws.Range("E8").Value = 9.4503885507905903
$ws.Range("F8").Value = 7.4224564551063859
$ws.Range("G8").Value = 6.1405538419242616
$ws.Range("H8").Value = 7.8362348089562586

# --- Row 9 ------------------------------------------------------------
$ws.Range("E9").Value = 12.05856261343326
$ws.Range("F9").Value = 12.869269925948988
$ws.Range("G9").Value = 12.224280084523084
$ws.Range("H9").Value = 12.015947868083865

# --- Row 10 -----------------------------------------------------------
$ws.Range("E10").Value = 12.539509799019589
$ws.Range("F10").Value = 11.669704858550739
$ws.Range("G10").Value = 10.487881438208545
$ws.Range("H10").Value = 9.0733515414336701

# --- Row 11 -----------------------------------------------------------
$ws.Range("E11").Value = 19.254146075966791
$ws.Range("F11").Value = 18.469923747895955
$ws.Range("G11").Value = 17.331663549619034
$ws.Range("H11").Value = 17.184771307810873

# --- Row 12 -----------------------------------------------------------
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0

# --- Row 13 -----------------------------------------------------------
$ws.Range("E13").Value = 14.586168576735787
$ws.Range("F13").Value = 15.129934195829577
$ws.Range("G13").Value = 14.464853778805054
$ws.Range("H13").Value = 13.872804190709424

# --- Row 14 -----------------------------------------------------------
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0

# --- Row 15 -----------------------------------------------------------
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# --- Row 16 -----------------------------------------------------------
$ws.Range("E16").Value = 16.074915790463564
$ws.Range("F16").Value = 15.146593404953883
$ws.Range("G16").Value = 14.375659010090962
$ws.Range("H16").Value = 13.935646874231191

# --- Row 17 -----------------------------------------------------------
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0

# --- Row 18 -----------------------------------------------------------
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 0
$ws.Range("H18").Value = 0

# --- Selection cursor, as saved in the source file --------------------
[void]$ws.Range("J5").Select()
